# Reproduce the "sheet2 now matches reader/sheet2" edit described in the
# commit message: Sheet2 - Numbers gains a 27th column (AA) of data,
# becomes the active/selected sheet (instead of Sheet4 - Dates), and
# Sheet4's page setup paper size is corrected.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Add the new AA column: AA1:AA30 = 100..129
for ($i = 0; $i -lt 30; $i++) {
    $ws2.Cells.Item($i + 1, 27).Value = 100 + $i
}

# Sheet2 becomes the active sheet/tab, with AA1:AA30 selected
$ws2.Activate()
$ws2.Range("AA1:AA30").Select()

# Sheet4's page setup: paper size goes from "automatic" (0) to Letter (9)
$ws4 = $wb.Worksheets.Item(4)
$ws4.PageSetup.PaperSize = 9
